$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column F to fit new data
$ws.Columns.Item(6).ColumnWidth = 10

# Add new trade row 7 - copy formatting (styles) down from row 6 first
$ws.Range("A6:I6").Copy() | Out-Null
$ws.Range("A7:I7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Now populate the new row's values
$ws.Range("A7").Value = 42649.644756944443
$ws.Range("B7").Value = $false
$ws.Range("C7").Value = 9842.75
$ws.Range("D7").Value = 9864.9500000000007
$ws.Range("E7").Value = 104.82
$ws.Range("F7").Value = 105.290001
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = 0.45
$ws.Range("I7").Value = $false
